$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.433768000000001
$ws.Range("H2").Value = 16.301304
$ws.Range("I2").Value = 0.1262505823713576
$ws.Range("J2").Value = 0.1262505823713576
$ws.Range("M2").Value = 9.084137666666667
$ws.Range("N2").Value = 27.252413
$ws.Range("O2").Value = 0.2765376761551382
$ws.Range("P2").Value = 0.2765376761551382
$ws.Range("Q2").Value = 49.36109656072801
$ws.Range("R2").Value = 444.2498690465521
$ws.Range("S2").Value = 0.0349130426622081
$ws.Range("T2").Value = 0.03491304266220809
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.433768000000001
$ws.Range("H3").Value = 16.301304
$ws.Range("I3").Value = 0.1262505823713576
$ws.Range("J3").Value = 0.1262505823713576
$ws.Range("O3").Value = 0.3707916163717078
$ws.Range("P3").Value = 0.3707916163717078
$ws.Range("Q3").Value = 66.18512541981602
$ws.Range("R3").Value = 595.6661287783442
$ws.Range("S3").Value = 0.04681265750534514
$ws.Range("T3").Value = 0.04681265750534513
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.433768000000001
$ws.Range("H4").Value = 16.301304
$ws.Range("I4").Value = 0.1262505823713576
$ws.Range("J4").Value = 0.1262505823713576
$ws.Range("M4").Value = 11.58507333333333
$ws.Range("N4").Value = 34.75522
$ws.Range("O4").Value = 0.3526707074731541
$ws.Range("P4").Value = 0.3526707074731541
$ws.Range("Q4").Value = 62.95060075632001
$ws.Range("R4").Value = 566.5554068068801
$ws.Range("S4").Value = 0.04452488220380441
$ws.Range("T4").Value = 0.04452488220380441
$ws.Range("G5").Value = 5.270503666666666
$ws.Range("I5").Value = 0.1224572262391479
$ws.Range("J5").Value = 0.1224572262391479
$ws.Range("M5").Value = 9.084137666666667
$ws.Range("N5").Value = 27.252413
$ws.Range("O5").Value = 0.2765376761551382
$ws.Range("P5").Value = 0.2765376761551382
$ws.Range("Q5").Value = 47.87798088067144
$ws.Range("R5").Value = 430.901827926043
$ws.Range("S5").Value = 0.03386403677257798
$ws.Range("T5").Value = 0.03386403677257797
$ws.Range("G6").Value = 5.270503666666666
$ws.Range("I6").Value = 0.1224572262391479
$ws.Range("J6").Value = 0.1224572262391479
$ws.Range("O6").Value = 0.3707916163717078
$ws.Range("P6").Value = 0.3707916163717078
$ws.Range("R6").Value = 577.7685973776211
$ws.Range("S6").Value = 0.04540611285360957
$ws.Range("T6").Value = 0.04540611285360956
$ws.Range("G7").Value = 5.270503666666666
$ws.Range("I7").Value = 0.1224572262391479
$ws.Range("J7").Value = 0.1224572262391479
$ws.Range("M7").Value = 11.58507333333333
$ws.Range("N7").Value = 34.75522
$ws.Range("O7").Value = 0.3526707074731541
$ws.Range("P7").Value = 0.3526707074731541
$ws.Range("Q7").Value = 61.05917148193556
$ws.Range("R7").Value = 549.53254333742
$ws.Range("S7").Value = 0.04318707661296039
$ws.Range("T7").Value = 0.04318707661296038
$ws.Range("G8").Value = 32.33527633333333
$ws.Range("H8").Value = 97.00582900000001
$ws.Range("I8").Value = 0.7512921913894945
$ws.Range("J8").Value = 0.7512921913894944
$ws.Range("M8").Value = 9.084137666666667
$ws.Range("N8").Value = 27.252413
$ws.Range("O8").Value = 0.2765376761551382
$ws.Range("P8").Value = 0.2765376761551382
$ws.Range("Q8").Value = 293.7381017017086
$ws.Range("R8").Value = 2643.642915315377
$ws.Range("S8").Value = 0.2077605967203521
$ws.Range("T8").Value = 0.2077605967203521
$ws.Range("G9").Value = 32.33527633333333
$ws.Range("H9").Value = 97.00582900000001
$ws.Range("I9").Value = 0.7512921913894945
$ws.Range("J9").Value = 0.7512921913894944
$ws.Range("O9").Value = 0.3707916163717078
$ws.Range("P9").Value = 0.3707916163717078
$ws.Range("Q9").Value = 393.8545627281244
$ws.Range("R9").Value = 3544.69106455312
$ws.Range("S9").Value = 0.2785728460127532
$ws.Range("T9").Value = 0.2785728460127531
$ws.Range("G10").Value = 32.33527633333333
$ws.Range("H10").Value = 97.00582900000001
$ws.Range("I10").Value = 0.7512921913894945
$ws.Range("J10").Value = 0.7512921913894944
$ws.Range("M10").Value = 11.58507333333333
$ws.Range("N10").Value = 34.75522
$ws.Range("O10").Value = 0.3526707074731541
$ws.Range("P10").Value = 0.3526707074731541
$ws.Range("Q10").Value = 374.6065475752645
$ws.Range("R10").Value = 3371.45892817738
$ws.Range("S10").Value = 0.2649587486563894
$ws.Range("T10").Value = 0.2649587486563893
